$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Tempo" sheet as the first sheet in the workbook.
# ---------------------------------------------------------------------------
$tempo = $wb.Worksheets.Add()
$tempo.Name = "Tempo"
$tempo.Move($wb.Worksheets.Item(1))

# ---------------------------------------------------------------------------
# 2. Populate "Tempo" with the sprint cost/time calculations.
# ---------------------------------------------------------------------------
$tempo.Range("B1").HorizontalAlignment = -4108

$tempo.Range("A2").Value = "Preço por hora"
$tempo.Range("B2").Value = 50
$tempo.Range("B2").HorizontalAlignment = -4108

$tempo.Range("A3").Value = "Quantidade de integrantes"
$tempo.Range("B3").Value = 7
$tempo.Range("B3").HorizontalAlignment = -4108

$tempo.Range("A4").Value = "Tempo diário"
$tempo.Range("B4").Value = 2.5
$tempo.Range("B4").HorizontalAlignment = -4108

$tempo.Range("A5").Value = "Quantidade de dias"
$tempo.Range("B5").Value = 21
$tempo.Range("B5").HorizontalAlignment = -4108

$tempo.Range("A6").Value = "Total horas sprint"
$tempo.Range("B6").Formula = "=B3*B4*B5"
$tempo.Range("B6").HorizontalAlignment = -4108
$tempo.Range("C6").Formula = "=B6/7"

$tempo.Range("A7").Value = "Total custo sprint"
$tempo.Range("B7").Formula = "=B2*B6"
$tempo.Range("B7").HorizontalAlignment = -4108
$tempo.Range("C7").Formula = "=B7/7"

$tempo.Range("A8").Value = "Total horas projeto"
$tempo.Range("B8").Formula = "=B6*4"
$tempo.Range("B8").HorizontalAlignment = -4108
$tempo.Range("C8").Formula = "=B8/7"

$tempo.Range("A9").Value = "Total custo projeto"
$tempo.Range("B9").Formula = "=B7*4"
$tempo.Range("B9").HorizontalAlignment = -4108
$tempo.Range("C9").Formula = "=B9/7"

$tempo.Columns.Item(1).ColumnWidth = 15.5

$tempo.PageSetup.LeftMargin = 0.511811024 * 72
$tempo.PageSetup.RightMargin = 0.511811024 * 72
$tempo.PageSetup.TopMargin = 0.78740157499999996 * 72
$tempo.PageSetup.BottomMargin = 0.78740157499999996 * 72
$tempo.PageSetup.HeaderMargin = 0.31496062000000002 * 72
$tempo.PageSetup.FooterMargin = 0.31496062000000002 * 72

$tempo.Range("B6").Select()

# ---------------------------------------------------------------------------
# 3. Update the "Sprint 0" sheet — registering the 1st sprint day.
# ---------------------------------------------------------------------------
$sprint0 = $wb.Worksheets.Item("Sprint 0")

# Mock-ups (row 6): the hours logged on day 5 (G6) are removed.
$sprint0.Range("G6").ClearContents()

# Dev. Front-End (row 7): day 1 (C7) becomes a computed 0.5+0.5 instead of a
# flat 5.
$sprint0.Range("C7").Formula = "=0.5+0.5"

# Dev. Back-End (row 8): the hours logged on day 9 (K8) are removed.
$sprint0.Range("K8").ClearContents()

# Testes (row 9): 0.25h logged on day 1 (C9).
$sprint0.Range("C9").Value = 0.25

# Banco de dados (row 10): hours logged on days 7 and 13 (I10, O10) removed.
$sprint0.Range("I10").ClearContents()
$sprint0.Range("O10").ClearContents()

# Outros (row 13): 1.5h logged on day 1 (C13).
$sprint0.Range("C13").Value = 1.5

# Estimate table: update the "horas estimadas por dificuldade" columns.
$sprint0.Range("C18").Value = 4
$sprint0.Range("C19").Value = 2
$sprint0.Range("C20").Value = 1

# "Horas estimadas" (B15) now pulls straight from the Tempo sheet total.
$sprint0.Range("B15").Formula = "=Tempo!B6"

$sprint0.Activate()
$sprint0.Range("E19").Select()
